$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to match the new longer company names
$ws.Columns.Item(1).ColumnWidth = 30.5 - 5/6

# Row 7: Asociatia Suntem Langa Tine
$ws.Range("A7").Value = "Asociatia Suntem Langa Tine"
$ws.Range("B7").Value = "suntemlangatine"
$ws.Range("C7").Value = "password"
$ws.Range("D7").Value = "suntem@gmail.com"
$ws.Range("E7").Value = 12345
$ws.Range("G7").Value = "ong"
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:suntem@gmail.com")
$ws.Range("D7").Style = "Hyperlink"

# Row 8: Asociatia Toti pentru Fericire
$ws.Range("A8").Value = "Asociatia Toti pentru Fericire"
$ws.Range("B8").Value = "fericiretoti"
$ws.Range("C8").Value = "password"
$ws.Range("D8").Value = "totifericire@gmail.com"
$ws.Range("E8").Value = 12345
$ws.Range("G8").Value = "ong"
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:totifericire@gmail.com")
$ws.Range("D8").Style = "Hyperlink"

# Row 9: Green everywhere
$ws.Range("A9").Value = "Green everywhere"
$ws.Range("B9").Value = "green"
$ws.Range("C9").Value = "password"
$ws.Range("D9").Value = "green@yahoo.com"
$ws.Range("E9").Value = 12345
$ws.Range("G9").Value = "ong"
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:green@yahoo.com")
$ws.Range("D9").Style = "Hyperlink"

# Update selection to match the author's final cursor position
[void]$ws.Range("A11").Select()
